$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-25 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-26 Monday", 2) | Out-Null
$d.Content.Find.Execute("5+38=43", $true, $false, $false, $false, $false, $true, 1, $false, "16+50=66", 2) | Out-Null
$d.Content.Find.Execute("77-67=10", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=61", 2) | Out-Null
$d.Content.Find.Execute("90+6=96", $true, $false, $false, $false, $false, $true, 1, $false, "78-8=70", 2) | Out-Null
$d.Content.Find.Execute("19+38=57", $true, $false, $false, $false, $false, $true, 1, $false, "65+2=67", 2) | Out-Null
$d.Content.Find.Execute("89-87=2", $true, $false, $false, $false, $false, $true, 1, $false, "81-45=36", 2) | Out-Null
$d.Content.Find.Execute("27+9=36", $true, $false, $false, $false, $false, $true, 1, $false, "94-70=24", 2) | Out-Null
$d.Content.Find.Execute("55-34=21", $true, $false, $false, $false, $false, $true, 1, $false, "14+8=22", 2) | Out-Null
$d.Content.Find.Execute("7+18=25", $true, $false, $false, $false, $false, $true, 1, $false, "32-12=20", 2) | Out-Null
$d.Content.Find.Execute("77-22=55", $true, $false, $false, $false, $false, $true, 1, $false, "97-24=73", 2) | Out-Null
$d.Content.Find.Execute("41+29=70", $true, $false, $false, $false, $false, $true, 1, $false, "47-15=32", 2) | Out-Null
$d.Content.Find.Execute("70+14=84", $true, $false, $false, $false, $false, $true, 1, $false, "72-6=66", 2) | Out-Null
$d.Content.Find.Execute("64-61=3", $true, $false, $false, $false, $false, $true, 1, $false, "85+12=97", 2) | Out-Null
$d.Content.Find.Execute("87-33=54", $true, $false, $false, $false, $false, $true, 1, $false, "35+21=56", 2) | Out-Null
$d.Content.Find.Execute("39-14=25", $true, $false, $false, $false, $false, $true, 1, $false, "66-51=15", 2) | Out-Null
$d.Content.Find.Execute("7+0=7", $true, $false, $false, $false, $false, $true, 1, $false, "96-41=55", 2) | Out-Null
$d.Content.Find.Execute("13+85=98", $true, $false, $false, $false, $false, $true, 1, $false, "9+50=59", 2) | Out-Null
$d.Content.Find.Execute("34+8=42", $true, $false, $false, $false, $false, $true, 1, $false, "83-56=27", 2) | Out-Null
$d.Content.Find.Execute("67+2=69", $true, $false, $false, $false, $false, $true, 1, $false, "73-66=7", 2) | Out-Null
$d.Content.Find.Execute("3-0=3", $true, $false, $false, $false, $false, $true, 1, $false, "61+36=97", 2) | Out-Null
$d.Content.Find.Execute("73-70=3", $true, $false, $false, $false, $false, $true, 1, $false, "68-59=9", 2) | Out-Null
$d.Content.Find.Execute("66+32=98", $true, $false, $false, $false, $false, $true, 1, $false, "61+20=81", 2) | Out-Null
$d.Content.Find.Execute("6+68=74", $true, $false, $false, $false, $false, $true, 1, $false, "98-89=9", 2) | Out-Null
$d.Content.Find.Execute("1+92=93", $true, $false, $false, $false, $false, $true, 1, $false, "19-11=8", 2) | Out-Null
$d.Content.Find.Execute("42+39=81", $true, $false, $false, $false, $false, $true, 1, $false, "24+66=90", 2) | Out-Null
$d.Content.Find.Execute("11+80=91", $true, $false, $false, $false, $false, $true, 1, $false, "90-17=73", 2) | Out-Null
$d.Content.Find.Execute("99-53=46", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=45", 2) | Out-Null
$d.Content.Find.Execute("81-64=17", $true, $false, $false, $false, $false, $true, 1, $false, "65+13=78", 2) | Out-Null
$d.Content.Find.Execute("68+15=83", $true, $false, $false, $false, $false, $true, 1, $false, "92-18=74", 2) | Out-Null
$d.Content.Find.Execute("39-23=16", $true, $false, $false, $false, $false, $true, 1, $false, "6+11=17", 2) | Out-Null
$d.Content.Find.Execute("57-29=28", $true, $false, $false, $false, $false, $true, 1, $false, "88-82=6", 2) | Out-Null
$d.Content.Find.Execute("66-30=36", $true, $false, $false, $false, $false, $true, 1, $false, "85-15=70", 2) | Out-Null
$d.Content.Find.Execute("68+12=80", $true, $false, $false, $false, $false, $true, 1, $false, "56+31=87", 2) | Out-Null
$d.Content.Find.Execute("61-7=54", $true, $false, $false, $false, $false, $true, 1, $false, "46+20=66", 2) | Out-Null
$d.Content.Find.Execute("70-66=4", $true, $false, $false, $false, $false, $true, 1, $false, "14+41=55", 2) | Out-Null
$d.Content.Find.Execute("40+13=53", $true, $false, $false, $false, $false, $true, 1, $false, "44+0=44", 2) | Out-Null
$d.Content.Find.Execute("98-59=39", $true, $false, $false, $false, $false, $true, 1, $false, "43+30=73", 2) | Out-Null
$d.Content.Find.Execute("7+38=45", $true, $false, $false, $false, $false, $true, 1, $false, "24+69=93", 2) | Out-Null
$d.Content.Find.Execute("60-9=51", $true, $false, $false, $false, $false, $true, 1, $false, "90-76=14", 2) | Out-Null
$d.Content.Find.Execute("14+60=74", $true, $false, $false, $false, $false, $true, 1, $false, "65+33=98", 2) | Out-Null
$d.Content.Find.Execute("78-72=6", $true, $false, $false, $false, $false, $true, 1, $false, "22+56=78", 2) | Out-Null
$d.Content.Find.Execute("12+9=21", $true, $false, $false, $false, $false, $true, 1, $false, "1+81=82", 2) | Out-Null
$d.Content.Find.Execute("63-32=31", $true, $false, $false, $false, $false, $true, 1, $false, "97-58=39", 2) | Out-Null
$d.Content.Find.Execute("50-3=47", $true, $false, $false, $false, $false, $true, 1, $false, "5+83=88", 2) | Out-Null
$d.Content.Find.Execute("55-10=45", $true, $false, $false, $false, $false, $true, 1, $false, "60-57=3", 2) | Out-Null
$d.Content.Find.Execute("42+15=57", $true, $false, $false, $false, $false, $true, 1, $false, "51-37=14", 2) | Out-Null
$d.Content.Find.Execute("37-34=3", $true, $false, $false, $false, $false, $true, 1, $false, "96-71=25", 2) | Out-Null
$d.Content.Find.Execute("58+20=78", $true, $false, $false, $false, $false, $true, 1, $false, "23-2=21", 2) | Out-Null
$d.Content.Find.Execute("89-82=7", $true, $false, $false, $false, $false, $true, 1, $false, "65-21=44", 2) | Out-Null
$d.Content.Find.Execute("84-3=81", $true, $false, $false, $false, $false, $true, 1, $false, "84-79=5", 2) | Out-Null
$d.Content.Find.Execute("4+92=96", $true, $false, $false, $false, $false, $true, 1, $false, "63-62=1", 2) | Out-Null
$d.Content.Find.Execute("96-35=61", $true, $false, $false, $false, $false, $true, 1, $false, "86+3=89", 2) | Out-Null
$d.Content.Find.Execute("89+6=95", $true, $false, $false, $false, $false, $true, 1, $false, "47+32=79", 2) | Out-Null
$d.Content.Find.Execute("56+5=61", $true, $false, $false, $false, $false, $true, 1, $false, "3+68=71", 2) | Out-Null
$d.Content.Find.Execute("72-45=27", $true, $false, $false, $false, $false, $true, 1, $false, "89-11=78", 2) | Out-Null
$d.Content.Find.Execute("72+10=82", $true, $false, $false, $false, $false, $true, 1, $false, "97-28=69", 2) | Out-Null
$d.Content.Find.Execute("67-10=57", $true, $false, $false, $false, $false, $true, 1, $false, "81-29=52", 2) | Out-Null
$d.Content.Find.Execute("68-56=12", $true, $false, $false, $false, $false, $true, 1, $false, "39+45=84", 2) | Out-Null
$d.Content.Find.Execute("68-45=23", $true, $false, $false, $false, $false, $true, 1, $false, "71-60=11", 2) | Out-Null
$d.Content.Find.Execute("14+79=93", $true, $false, $false, $false, $false, $true, 1, $false, "68+22=90", 2) | Out-Null
$d.Content.Find.Execute("48+24=72", $true, $false, $false, $false, $false, $true, 1, $false, "12+12=24", 2) | Out-Null
$d.Content.Find.Execute("18+7=25", $true, $false, $false, $false, $false, $true, 1, $false, "29+69=98", 2) | Out-Null
$d.Content.Find.Execute("64-19=45", $true, $false, $false, $false, $false, $true, 1, $false, "22-7=15", 2) | Out-Null
$d.Content.Find.Execute("76+3=79", $true, $false, $false, $false, $false, $true, 1, $false, "32-29=3", 2) | Out-Null
$d.Content.Find.Execute("69-6=63", $true, $false, $false, $false, $false, $true, 1, $false, "97-58=39", 2) | Out-Null
$d.Content.Find.Execute("27-26=1", $true, $false, $false, $false, $false, $true, 1, $false, "18+8=26", 2) | Out-Null
$d.Content.Find.Execute("74-2=72", $true, $false, $false, $false, $false, $true, 1, $false, "22+8=30", 2) | Out-Null
$d.Content.Find.Execute("86-49=37", $true, $false, $false, $false, $false, $true, 1, $false, "81-26=55", 2) | Out-Null
$d.Content.Find.Execute("90-54=36", $true, $false, $false, $false, $false, $true, 1, $false, "43-12=31", 2) | Out-Null
$d.Content.Find.Execute("3+40=43", $true, $false, $false, $false, $false, $true, 1, $false, "60-14=46", 2) | Out-Null
$d.Content.Find.Execute("18+19=37", $true, $false, $false, $false, $false, $true, 1, $false, "0+41=41", 2) | Out-Null
$d.Content.Find.Execute("44+32=76", $true, $false, $false, $false, $false, $true, 1, $false, "38+29=67", 2) | Out-Null
$d.Content.Find.Execute("59+36=95", $true, $false, $false, $false, $false, $true, 1, $false, "55-19=36", 2) | Out-Null
$d.Content.Find.Execute("38-12=26", $true, $false, $false, $false, $false, $true, 1, $false, "92-53=39", 2) | Out-Null
$d.Content.Find.Execute("96-86=10", $true, $false, $false, $false, $false, $true, 1, $false, "71-25=46", 2) | Out-Null
$d.Content.Find.Execute("30-23=7", $true, $false, $false, $false, $false, $true, 1, $false, "89-18=71", 2) | Out-Null
$d.Content.Find.Execute("33-20=13", $true, $false, $false, $false, $false, $true, 1, $false, "63-63=0", 2) | Out-Null
$d.Content.Find.Execute("2+15=17", $true, $false, $false, $false, $false, $true, 1, $false, "30+10=40", 2) | Out-Null
$d.Content.Find.Execute("26+68=94", $true, $false, $false, $false, $false, $true, 1, $false, "27-23=4", 2) | Out-Null
$d.Content.Find.Execute("37-30=7", $true, $false, $false, $false, $false, $true, 1, $false, "99-83=16", 2) | Out-Null
$d.Content.Find.Execute("0+33=33", $true, $false, $false, $false, $false, $true, 1, $false, "16+42=58", 2) | Out-Null
$d.Content.Find.Execute("22+6=28", $true, $false, $false, $false, $false, $true, 1, $false, "97-7=90", 2) | Out-Null
$d.Content.Find.Execute("11+44=55", $true, $false, $false, $false, $false, $true, 1, $false, "67-32=35", 2) | Out-Null
$d.Content.Find.Execute("55+36=91", $true, $false, $false, $false, $false, $true, 1, $false, "15+5=20", 2) | Out-Null
$d.Content.Find.Execute("56+7=63", $true, $false, $false, $false, $false, $true, 1, $false, "9+0=9", 2) | Out-Null
$d.Content.Find.Execute("44-14=30", $true, $false, $false, $false, $false, $true, 1, $false, "13+38=51", 2) | Out-Null
$d.Content.Find.Execute("23+15=38", $true, $false, $false, $false, $false, $true, 1, $false, "62+3=65", 2) | Out-Null
$d.Content.Find.Execute("76-54=22", $true, $false, $false, $false, $false, $true, 1, $false, "54+1=55", 2) | Out-Null
$d.Content.Find.Execute("84+7=91", $true, $false, $false, $false, $false, $true, 1, $false, "60+31=91", 2) | Out-Null
$d.Content.Find.Execute("63+23=86", $true, $false, $false, $false, $false, $true, 1, $false, "1+82=83", 2) | Out-Null
$d.Content.Find.Execute("19+52=71", $true, $false, $false, $false, $false, $true, 1, $false, "70+11=81", 2) | Out-Null
$d.Content.Find.Execute("88-74=14", $true, $false, $false, $false, $false, $true, 1, $false, "16+22=38", 2) | Out-Null
$d.Content.Find.Execute("66+11=77", $true, $false, $false, $false, $false, $true, 1, $false, "58+22=80", 2) | Out-Null
$d.Content.Find.Execute("56-8=48", $true, $false, $false, $false, $false, $true, 1, $false, "0+82=82", 2) | Out-Null
$d.Content.Find.Execute("29+43=72", $true, $false, $false, $false, $false, $true, 1, $false, "28+52=80", 2) | Out-Null
$d.Content.Find.Execute("32+38=70", $true, $false, $false, $false, $false, $true, 1, $false, "19+45=64", 2) | Out-Null
$d.Content.Find.Execute("60-5=55", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=91", 2) | Out-Null
$d.Content.Find.Execute("58+4=62", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=95", 2) | Out-Null
$d.Content.Find.Execute("92-22=70", $true, $false, $false, $false, $false, $true, 1, $false, "4+3=7", 2) | Out-Null
$d.Content.Find.Execute("95-77=18", $true, $false, $false, $false, $false, $true, 1, $false, "5+45=50", 2) | Out-Null
$d.Content.Find.Execute("29-15=14", $true, $false, $false, $false, $false, $true, 1, $false, "76-53=23", 2) | Out-Null
